{"js": "// Remove the stray \".\" that directly follows the \"{{ \u0161k_sat_z }}\" merge\n// field (it currently reads \"{{ \u0161k_sat_z }}. \u0161kolski sat.\" and should read\n// \"{{ \u0161k_sat_z }} \u0161kolski sat.\").\nconst body = context.document.body;\n\n// The merge field is split across several runs (\"{{ \u0161k_sat\" + \"_z\" + \" }}\"),\n// but Word.js search matches across run boundaries inside one paragraph, so\n// searching for the whole tag plus the trailing dot finds it in one hit.\nconst hits = body.search(\"{{ \u0161k_sat_z }}.\", { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"{{ \u0161k_sat_z }}.\" in the document body.');\n}\n\nfor (const hit of hits.items) {\n  // Narrow the hit down to just the trailing \".\" character so only that\n  // run's text is touched, then delete it outright (removing the run,\n  // exactly like the authored edit).\n  const dot = hit.search(\".\", { matchCase: true, matchWildcards: false });\n  dot.load(\"text\");\n  await context.sync();\n\n  const last = dot.items[dot.items.length - 1];\n  last.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the stray \".\" that directly follows the \"{{ \u0161k_sat_z }}\" merge\n# field (it currently reads \"{{ \u0161k_sat_z }}. \u0161kolski sat.\" and should read\n# \"{{ \u0161k_sat_z }} \u0161kolski sat.\").\n\n$d = $word.ActiveDocument\n\n# Locate the merge field together with its trailing dot. Find matches across\n# the run boundaries that split up \"{{ \u0161k_sat\" + \"_z\" + \" }}\" because it\n# operates on the story's plain text, not on individual runs.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Forward = $true\n$find.Wrap = 0            # wdFindStop - do not wrap around, there is a single hit\n$find.MatchCase = $true\n$find.MatchWildcards = $false\n$find.Text = \"{{ \u0161k_sat_z }}.\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"{{ \u0161k_sat_z }}.\" in the document.'\n}\n\n# $rng now spans the matched text \"{{ \u0161k_sat_z }}.\"; narrow it down to just\n# the final \".\" character so only that run is touched, then delete it - the\n# same net effect as removing the whole <w:r> that holds the lone \".\".\n$rng.Start = $rng.End - 1\n$rng.Delete()\n"}
